$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Row 29 / Row 30: coin identity swapped, with refreshed price/volume ---
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.597.54'
$ws.Range("E29").Value = '  +3.02%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D30" '1.02'
$ws.Range("E30").Value = '  +1.90%  '

# --- Price / Volume(1h) refresh for remaining rows ---
$ws.Range("D2").Value = '63.347.04'
$ws.Range("E2").Value = '  +2.52%  '
$ws.Range("D3").Value = '2.475.53'
$ws.Range("E3").Value = '  +2.59%  '
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.32%  '
Set-TextValue "D5" '578.35'
$ws.Range("E5").Value = '  +1.74%  '
Set-TextValue "D6" '147.09'
$ws.Range("E6").Value = '  +2.13%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = '2.472.31'
$ws.Range("E9").Value = '  +1.84%  '
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("E11").Value = '  +1.41%  '
Set-TextValue "D12" '5.30'
$ws.Range("E12").Value = '  +1.37%  '
Set-TextValue "D13" '0.356'
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("E14").Value = '  +9.62%  '
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '2.918.52'
$ws.Range("E16").Value = '  +2.92%  '
$ws.Range("D17").Value = '63.200.25'
$ws.Range("E17").Value = '  +2.45%  '
$ws.Range("D18").Value = '2.467.29'
$ws.Range("E18").Value = '  +1.87%  '
Set-TextValue "D19" '7.94'
$ws.Range("E19").Value = '  -0.07%  '
Set-TextValue "D20" '11.13'
$ws.Range("E20").Value = '  +3.65%  '
Set-TextValue "D21" '331.27'
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("E22").Value = '  +10.47%  '
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("E24").Value = '  +0.07%  '
Set-TextValue "D25" '66.56'
$ws.Range("E25").Value = '  +2.11%  '
Set-TextValue "D26" '676.16'
$ws.Range("E26").Value = '  +8.86%  '
Set-TextValue "D27" '9.19'
$ws.Range("E27").Value = '  +9.47%  '
$ws.Range("E28").Value = '  +6.03%  '
$ws.Range("E31").Value = '  +3.81%  '
Set-TextValue "D32" '8.21'
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("E33").Value = '  +4.24%  '
Set-TextValue "D34" '0.139'
$ws.Range("E34").Value = '  +3.33%  '
Set-TextValue "D35" '1.57'
$ws.Range("E35").Value = '  +5.63%  '
Set-TextValue "D36" '0.998'
$ws.Range("E36").Value = '  +0.04%  '
Set-TextValue "D37" '4.82'
$ws.Range("E37").Value = '  +4.17%  '
Set-TextValue "D38" '5.59'
$ws.Range("E38").Value = '  +4.40%  '
$ws.Range("E39").Value = '  +1.08%  '
Set-TextValue "D40" '153.64'
$ws.Range("E40").Value = '  +0.26%  '
Set-TextValue "D41" '18.91'
$ws.Range("E41").Value = '  +2.56%  '
Set-TextValue "D42" '2.79'
$ws.Range("E42").Value = '  +7.50%  '
$ws.Range("E43").Value = '  +3.95%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").Value = '0.0₆0301'
$ws.Range("E46").Value = '  +5.79%  '
Set-TextValue "D47" '15.15'
$ws.Range("E47").Value = '  +27.77%  '
Set-TextValue "D48" '147.12'
$ws.Range("E48").Value = '  +3.35%  '
Set-TextValue "D49" '3.65'
$ws.Range("E49").Value = '  +2.64%  '
Set-TextValue "D50" '20.90'
$ws.Range("E50").Value = '  +4.49%  '
Set-TextValue "D51" '0.609'
$ws.Range("E51").Value = '  +1.96%  '
